$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95, shifting rows 95:125 down to 96:126.
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new data point.
$ws.Range("A95").Value = 6
$ws.Range("B95").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C95").Value = 'Metropolitana'
$ws.Range("D95").Value = 44524
$ws.Range("D95").NumberFormat = $ws.Range("D94").NumberFormat
$ws.Range("E95").Value = 13
$ws.Range("F95").Value = 100112029
$ws.Range("G95").Value = 'Orégano'
$ws.Range("H95").Value = 'Sin especificar'
$ws.Range("I95").Value = 'Primera'
$ws.Range("J95").Value = 34
$ws.Range("K95").Value = 8500
$ws.Range("L95").Value = 9000
$ws.Range("M95").Value = 8721
$ws.Range("N95").Value = '$/docena de atados'
$ws.Range("O95").Value = 'Región Metropolitana'
$ws.Range("P95").Value = 2907
$ws.Range("Q95").Value = 3
$ws.Range("R95").Value = 'Hortaliza'
